# The commit swaps the presentation's two embedded themes: the theme used
# by the slide master (ppt/theme/theme1.xml, originally the "Integral"
# palette) is replaced with what used to be the notes master's theme
# (ppt/theme/theme2.xml, the stock "Office Theme" palette). The font
# scheme / format scheme are identical between the two themes already, so
# nothing else needs to change.
#
# PowerPoint's object model exposes the live theme colour slots for the
# deck's (single) slide master through any Slide's ThemeColorScheme
# property (Colors 1-12 == dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) - editing them rewrites ppt/theme/theme1.xml directly, the
# same effect a user gets by picking a different built-in colour theme
# from the Design tab. RGB values below are packed as 0xBBGGRR (the
# native OLE COLORREF layout PowerPoint's ColorFormat.RGB uses).

$p = $ppt.ActivePresentation
$cs = $p.Slides.Item(1).ThemeColorScheme

$cs.Item(1).RGB  = 0x000000   # dk1      000000
$cs.Item(2).RGB  = 0xFFFFFF   # lt1      FFFFFF
$cs.Item(3).RGB  = 0x6A5444   # dk2      44546A
$cs.Item(4).RGB  = 0xE6E6E7   # lt2      E7E6E6
$cs.Item(5).RGB  = 0xD59B5B   # accent1  5B9BD5
$cs.Item(6).RGB  = 0x317DED   # accent2  ED7D31
$cs.Item(7).RGB  = 0xA5A5A5   # accent3  A5A5A5
$cs.Item(8).RGB  = 0x00C0FF   # accent4  FFC000
$cs.Item(9).RGB  = 0xC47244   # accent5  4472C4
$cs.Item(10).RGB = 0x47AD70   # accent6  70AD47
$cs.Item(11).RGB = 0xC16305   # hlink    0563C1
$cs.Item(12).RGB = 0x724F95   # folHlink 954F72
